# Collated Test Set Results.xlsx — "Add opus-big AoN wce results"
#
# - opus_base Test  (sheet1 / codeName Sheet2): updates B7/B8/B11/B12 labels
#   and fills in two brand-new result rows (16 & 17: "1.0 train sampled
#   glossary, bands 5, weight 1.25" / "...proportion 0.6, weight 1.5").
# - opus_big Test   (sheet2 / codeName Sheet3): updates B7/B8/B11/B12 labels
#   and fills in one brand-new result row (16: "1.0 train sampled glossary,
#   5 bands, weight 1.75").
# - Active sheet/selection flips from "opus_big Test " (B35-ish) to
#   "opus_base Test" (B7).

$wb = $excel.ActiveWorkbook

$wsBase = $wb.Worksheets.Item("opus_base Test")
$wsBig  = $wb.Worksheets.Item("opus_big Test ")

# ---------------------------------------------------------------------
# opus_base Test
# ---------------------------------------------------------------------

# Relabel existing rows (values for these rows stay the same; only the
# "Brief Description" text in column B changes).
$wsBase.Cells.Item(7, 2).Value  = "0.1 train sampled glossary, weight 1.25"
$wsBase.Cells.Item(8, 2).Value  = "0.1 train sample, full glossary, upper weight 1.25, bands 6"
$wsBase.Cells.Item(11, 2).Value = "0.25 train sampled glossary, weight 1.5"
$wsBase.Cells.Item(12, 2).Value = "0.25 train sample, full glossary, upper weight 1.5, bands 6"

# Row 16 ("Simple Adaptive WCE" section) — new result: sampled-glossary,
# bands 5, weight 1.25. B16 needs to pick up the same format as B17
# (fillId/border combo, xf 40) since it was previously a blank styled cell.
$wsBase.Range("B17").Copy() | Out-Null
$wsBase.Range("B16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsBase.Cells.Item(16, 2).Value  = "1.0 train sampled glossary, bands 5, weight 1.25"
$wsBase.Cells.Item(16, 3).Value  = 42.950499999999998
$wsBase.Cells.Item(16, 4).Value  = 33.340000000000003
$wsBase.Cells.Item(16, 5).Value  = 52.65
$wsBase.Cells.Item(16, 6).Value  = 69.540000000000006
$wsBase.Cells.Item(16, 7).Value  = 0.48480000000000001
$wsBase.Cells.Item(16, 8).Value  = 0.68620000000000003
$wsBase.Cells.Item(16, 9).Value  = 960
$wsBase.Cells.Item(16, 10).Value = 0.68979999999999997
$wsBase.Cells.Item(16, 11).Value = 965
$wsBase.Cells.Item(16, 12).Value = 0.74909999999999999
$wsBase.Cells.Item(16, 13).Value = 1048

# Row 17 ("All-or-Nothing Adaptive WCE" section) — new result: sampled
# glossary, proportion 0.6, weight 1.5. Formatting here is already in
# place (xf 40/41), only values are required.
$wsBase.Cells.Item(17, 2).Value  = "1.0 train sampled glossary, proportion 0.6, weight 1.5"
$wsBase.Cells.Item(17, 3).Value  = 43.082599999999999
$wsBase.Cells.Item(17, 4).Value  = 33.409999999999997
$wsBase.Cells.Item(17, 5).Value  = 52.72
$wsBase.Cells.Item(17, 6).Value  = 69.44
$wsBase.Cells.Item(17, 7).Value  = 0.48599999999999999
$wsBase.Cells.Item(17, 8).Value  = 0.68189999999999995
$wsBase.Cells.Item(17, 9).Value  = 954
$wsBase.Cells.Item(17, 10).Value = 0.68479999999999996
$wsBase.Cells.Item(17, 11).Value = 958
$wsBase.Cells.Item(17, 12).Value = 0.74619999999999997
$wsBase.Cells.Item(17, 13).Value = 1044

# ---------------------------------------------------------------------
# opus_big Test
# ---------------------------------------------------------------------

$wsBig.Cells.Item(7, 2).Value  = "0.1 train sampled glossary, weight 1.25"
$wsBig.Cells.Item(8, 2).Value  = "0.1 train sample, full glossary, upper weight 1.25, bands 6"
$wsBig.Cells.Item(11, 2).Value = "0.75 train sampled glossary, weight 1.25"
$wsBig.Cells.Item(12, 2).Value = "0.75 train sample, full glossary, upper weight 1.25, bands 6"

# Row 16 ("Simple Adaptive WCE" section) — new opus-big AoN WCE result:
# sampled glossary, 5 bands, weight 1.75. Existing blank-cell formatting
# (xf 30/37) already matches the target, only values are required.
$wsBig.Cells.Item(16, 2).Value  = "1.0 train sampled glossary, 5 bands, weight 1.75"
$wsBig.Cells.Item(16, 3).Value  = 43.872199999999999
$wsBig.Cells.Item(16, 4).Value  = 34.24
$wsBig.Cells.Item(16, 5).Value  = 53.31
$wsBig.Cells.Item(16, 6).Value  = 67.900000000000006
$wsBig.Cells.Item(16, 7).Value  = 0.51749999999999996
$wsBig.Cells.Item(16, 8).Value  = 0.74409999999999998
$wsBig.Cells.Item(16, 9).Value  = 1041
$wsBig.Cells.Item(16, 10).Value = 0.747
$wsBig.Cells.Item(16, 11).Value = 1045
$wsBig.Cells.Item(16, 12).Value = 0.747
$wsBig.Cells.Item(16, 13).Value = 1045

# ---------------------------------------------------------------------
# Active sheet / selection — "opus_base Test" becomes the active tab
# with B7 selected; "opus_big Test " keeps a plain (non-active) B35
# selection.
# ---------------------------------------------------------------------

$wsBig.Activate()
$wsBig.Range("B35").Select() | Out-Null

$wsBase.Activate()
$wsBase.Range("B7").Select() | Out-Null

Write-Output "done"
